# Daily attendance processing - 2025-12-04 07:28:57
#
# Re-orders the "Recorded By" entries (column G) so that, for any cell whose
# value is a comma-separated list that does NOT contain "admin@admin.com",
# the last two entries in the list are swapped. Cells with a single entry,
# or that reference admin@admin.com, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recordedByCol = 7  # Column G = "Recorded By"
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $value = $ws.Cells.Item($r, $recordedByCol).Value()
    if ([string]::IsNullOrEmpty($value)) {
        continue
    }
    if ($value.Contains("admin@admin.com")) {
        continue
    }

    $parts = $value -split ",\s*"
    if ($parts.Count -ge 2) {
        $last = $parts.Count - 1
        $secondLast = $parts.Count - 2

        $tmp = $parts[$last]
        $parts[$last] = $parts[$secondLast]
        $parts[$secondLast] = $tmp

        $newValue = [string]::Join(", ", $parts)
        $ws.Cells.Item($r, $recordedByCol).Value = $newValue
    }
}
